# Update cryptos list (GitHub Actions style refresh of prices / 1h volume %)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "price" value into column D while forcing it to remain
# plain text (many of the scraped price strings look numeric, e.g. "1.010",
# and would otherwise be reinterpreted as numbers and lose trailing zeros).
function Set-PriceText($row, $text) {
    $ws.Range("D$row").Value = "'" + $text
}

# Rows whose Coin/Link/Price/Volume just refresh in place -----------------

Set-PriceText 2 "27.005.04"
$ws.Range("E2").Value = "  -0.50%  "

Set-PriceText 3 "1.829.18"
$ws.Range("E3").Value = "  +0.16%  "

Set-PriceText 4 "1.010"
$ws.Range("E4").Value = "  -0.18%  "

Set-PriceText 5 "311.81"
$ws.Range("E5").Value = "  -0.37%  "

Set-PriceText 6 "1.008"
$ws.Range("E6").Value = "  -0.19%  "

Set-PriceText 7 "0.4648"
$ws.Range("E7").Value = "  -1.14%  "

Set-PriceText 8 "0.3709"
$ws.Range("E8").Value = "  +1.71%  "

Set-PriceText 9 "0.07370"
$ws.Range("E9").Value = "  -0.40%  "

Set-PriceText 10 "0.8769"
$ws.Range("E10").Value = "  -0.35%  "

# Rows 11 and 12 swap places: TRON/Solana order flips ----------------------

$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-PriceText 11 "19.99"
$ws.Range("E11").Value = "  -1.82%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-PriceText 12 "0.07874"
$ws.Range("E12").Value = "  +7.36%  "

Set-PriceText 13 "1.777.44"
$ws.Range("E13").Value = "  -8.46%  "

Set-PriceText 14 "5.361"
$ws.Range("E14").Value = "  -0.39%  "

Set-PriceText 15 "6.588"
$ws.Range("E15").Value = "  +1.15%  "

Set-PriceText 16 "92.07"
$ws.Range("E16").Value = "  -1.20%  "

Set-PriceText 17 "1.010"
$ws.Range("E17").Value = "  +0.17%  "

Set-PriceText 18 "0.000008889"
$ws.Range("E18").Value = "  +2.05%  "

Set-PriceText 19 "1.008"
$ws.Range("E19").Value = "  -0.30%  "

# Row 20 (Avalanche): price unchanged, only the 1h volume % moves
$ws.Range("E20").Value = "  +0.51%  "

Set-PriceText 21 "26.977.08"
$ws.Range("E21").Value = "  -2.64%  "

Set-PriceText 22 "5.166"
$ws.Range("E22").Value = "  -1.45%  "

Set-PriceText 23 "10.59"
$ws.Range("E23").Value = "  +0.17%  "

Set-PriceText 24 "1.985.45"
$ws.Range("E24").Value = "  -5.71%  "

Set-PriceText 25 "152.61"
$ws.Range("E25").Value = "  +0.68%  "

Set-PriceText 26 "1.834"
$ws.Range("E26").Value = "  -2.56%  "

Set-PriceText 27 "18.28"
$ws.Range("E27").Value = "  -1.45%  "

Set-PriceText 28 "2.101"
$ws.Range("E28").Value = "  -1.55%  "

Set-PriceText 29 "5.131"
$ws.Range("E29").Value = "  -0.77%  "

Set-PriceText 30 "115.61"
$ws.Range("E30").Value = "  -0.50%  "

Set-PriceText 31 "0.08884"
$ws.Range("E31").Value = "  -0.60%  "

Set-PriceText 32 "2.984"
$ws.Range("E32").Value = "  +1.46%  "

# Rows 33 and 34 swap places: ImmutableX/Filecoin order flips --------------

$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-PriceText 33 "4.451"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-PriceText 34 "0.7283"
$ws.Range("E34").Value = "  -1.71%  "

# Row 35 (ARBITRUM): price unchanged, only the 1h volume % moves
$ws.Range("E35").Value = "  -2.62%  "

# Row 36 (RenderToken): price unchanged, only the 1h volume % moves
$ws.Range("E36").Value = "  +3.41%  "

Set-PriceText 37 "1.080"
$ws.Range("E37").Value = "  -0.81%  "

Set-PriceText 38 "0.01956"
$ws.Range("E38").Value = "  +0.55%  "

Set-PriceText 39 "0.05248"
$ws.Range("E39").Value = "  -0.86%  "

Set-PriceText 40 "7.323"
$ws.Range("E40").Value = "  +1.90%  "

Set-PriceText 41 "2.932"
$ws.Range("E41").Value = "  -0.38%  "

Set-PriceText 42 "0.5200"
$ws.Range("E42").Value = "  -1.04%  "

Set-PriceText 43 "0.1628"
$ws.Range("E43").Value = "  -0.80%  "

Set-PriceText 44 "0.8602"
$ws.Range("E44").Value = "  -14.89%  "

Set-PriceText 45 "8.225"
$ws.Range("E45").Value = "  -1.70%  "

Set-PriceText 46 "0.4849"
$ws.Range("E46").Value = "  -0.58%  "

Set-PriceText 47 "1.009"
$ws.Range("E47").Value = "  -0.22%  "

# Row 48 (EnergySwap): price unchanged, only the 1h volume % moves
$ws.Range("E48").Value = "  -0.97%  "

Set-PriceText 49 "102.83"
$ws.Range("E49").Value = "  -1.48%  "

# Row 50 (NEARProtocol): price unchanged, only the 1h volume % moves
$ws.Range("E50").Value = "  -1.38%  "

Set-PriceText 51 "0.06237"
$ws.Range("E51").Value = "  -1.02%  "
